$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ("ECs", "Col2a1", "Mag", "ECs", 3, 1, 0.01944166666666667, 0.058325, 0.085855768197522, 0.085855768197522, 1, 0.3333333333333333, 0.09045733333333333, 0.271372, 0.0574491187297735, 0.0574491187297735, 0.001758641322222222, 0.0158277719, 0.004932338220815353, 0.004932338220815353),
    ("ECs", "Col2a1", "Mag", "sCs", 3, 1, 0.01944166666666667, 0.058325, 0.085855768197522, 0.085855768197522, 3, 1, 1.484107, 4.452321, 0.9425508812702265, 0.9425508812702265, 0.02885351359166667, 0.259681622325, 0.08092342997670664, 0.08092342997670664),
    ("FAPs", "Col2a1", "Mag", "ECs", 1, 0.3333333333333333, 0.2021233333333333, 0.60637, 0.8925908643280139, 0.8925908643280139, 1, 0.3333333333333333, 0.09045733333333333, 0.271372, 0.0574491187297735, 0.0574491187297735, 0.01828353773777778, 0.16455183964, 0.05127855854189122, 0.05127855854189122),
    ("FAPs", "Col2a1", "Mag", "sCs", 1, 0.3333333333333333, 0.2021233333333333, 0.60637, 0.8925908643280139, 0.8925908643280139, 3, 1, 1.484107, 4.452321, 0.9425508812702265, 0.9425508812702265, 0.2999726538633333, 2.69975388477, 0.8413123057861227, 0.8413123057861227),
    ("sCs", "Col2a1", "Mag", "ECs", 2, 0.6666666666666666, 0.004880666666666667, 0.014642, 0.02155336747446408, 0.02155336747446408, 1, 0.3333333333333333, 0.09045733333333333, 0.271372, 0.0574491187297735, 0.0574491187297735, 0.0004414920915555556, 0.003973428824, 0.001238221967066925, 0.001238221967066925),
    ("sCs", "Col2a1", "Mag", "sCs", 2, 0.6666666666666666, 0.004880666666666667, 0.014642, 0.02155336747446408, 0.02155336747446408, 3, 1, 1.484107, 4.452321, 0.9425508812702265, 0.9425508812702265, 0.007243431564666667, 0.06519088408200001, 0.02031514550739715, 0.02031514550739715),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $vals[$j]
    }
}
